# Adds the Ankara hospital rows (67 rows) to Sayfa1, grouped by sub-category
# with a blank separator row between each group, mirroring the pattern
# already used throughout the sheet for other provinces.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$il = "Ankara"

# Each group: (Tur, Hastane Adi list). A blank row separates each group,
# matching the existing layout convention used for every other province.
$groups = @(
    @{ Tur = "Kamu"; Names = @(
        "Ankara Bilkent Şehir Hastanesi",
        "Ankara Etlik Şehir Hastanesi"
    ) },
    @{ Tur = "Kamu"; Names = @(
        "Ankara Atatürk Sanatoryum Eğitim ve Araştırma Hastanesi",
        "Ankara Eğitim ve Araştırma Hastanesi",
        "Gülhane Eğitim ve Araştırma Hastanesi",
        "Dışkapı Yıldırım Beyazıt Eğitim ve Araştırma Hastanesi",
        "Dr. Abdurrahman Yurtaslan Ankara Onkoloji Eğitim ve Araştırma Hastanesi",
        "Etlik Zübeyde Hanım Kadın Hastalıkları Eğitim ve Araştırma Hastanesi",
        "Sincan Eğitim ve Araştırma Hastanesi",
        "Ankara Dr. Sami Ulus Çocuk Sağlığı ve Hastalıkları Eğitim ve Araştırma Hastanesi",
        "Ankara Gaziler Fizik Tedavi ve Rehabilitasyon Eğitim ve Araştırma Hastanesi",
        "Ulucanlar Göz Eğitim ve Araştırma Hastanesi",
        "Yenimahalle Eğitim ve Araştırma Hastanesi"
    ) },
    @{ Tur = "Kamu"; Names = @(
        "75. Yıl Ağız ve Diş Sağlığı Hastanesi",
        "Osmanlı Ağız ve Diş Sağlığı Hastanesi",
        "Polatlı Ağız ve Diş Sağlığı Hastanesi",
        "Balgat Ağız ve Diş Sağlığı Hastanesi",
        "Etimesgut Ağız ve Diş Sağlığı Hastanesi",
        "Ülkü Ulusoy Ağız ve Diş Sağlığı Hastanesi",
        "Karapürçek Ağız ve Diş Sağlığı Hastanesi",
        "Mamak Ağız ve Diş Sağlığı Hastanesi",
        "Pursaklar Ağız ve Diş Sağlığı Hastanesi",
        "Sincan Ağız ve Diş Sağlığı Hastanesi",
        "Tepebaşı Ağız ve Diş Sağlığı Hastanesi",
        "Topraklık Ağız ve Diş Sağlığı Hastanesi"
    ) },
    @{ Tur = "Kamu"; Names = @(
        "Ankara Beştepe Devlet Hastanesi",
        "Akyurt Devlet Hastanesi",
        "Ayaş Şehit Mehmet Çifci Devlet Hastanesi",
        "Bala Devlet Hastanesi",
        "Beştepe Devlet Hastanesi",
        "Beypazarı Devlet Hastanesi",
        "Beytepe Şehit Murat Erdi Eker Devlet Hastanesi",
        "Ceza İnfaz Kurumları Kampüs Devlet Hastanesi",
        "Çubuk Halil Şıvgın Devlet Hastanesi",
        "Elmadağ Dr. Hulusi Alataş Devlet Hastanesi",
        "Etimesgut Şehit Sait Ertürk Devlet Hastanesi",
        "Gazi Mustafa Kemal Mesleki ve Çevresel Hastalıklar Hastanesi",
        "Gölbaşı Şehit Ahmet Özsoy Devlet Hastanesi",
        "Güdül İlçe Entegre Devlet Hastanesi",
        "Haymana Devlet Hastanesi",
        "Kahramankazan Devlet Hastanesi",
        "Kalecik İlçe Entegre Devlet Hastanesi",
        "Kızılcahamam Devlet Hastanesi",
        "Mamak Devlet Hastanesi",
        "Nallıhan Devlet Hastanesi",
        "Polatlı Duatepe Devlet Hastanesi",
        "Pursaklar Devlet Hastanesi",
        "Şereflikoçhisar Devlet Hastanesi",
        "29 Mayıs Devlet Hastanesi"
    ) },
    @{ Tur = "Üniversite"; Names = @(
        "Ankara Üniversitesi Cebeci Araştırma ve Uygulama Hastanesi",
        "Ankara Üniversitesi İbni Sina Araştırma ve Uygulama Hastanesi",
        "Başkent Üniversitesi Ankara Hastanesi",
        "Gazi Üniversitesi Tıp Fakültesi Hastanesi",
        "Hacettepe Üniversitesi Beytepe Gün Hastanesi",
        "Hacettepe Üniversitesi Erişkin Hastanesi",
        "Hacettepe Üniversitesi İhsan Doğramacı Çocuk Hastanesi",
        "Hacettepe Üniversitesi Onkoloji Hastanesi",
        "Ufuk Üniversitesi Dr. Rıdvan Ege Hastanesi"
    ) },
    @{ Tur = "Özel"; Names = @(
        "Bayındır Hastanesi",
        "Dünya Göz Hastanesi (Çankaya)",
        "Dünya Göz Hastanesi Keçiören",
        "Dünya Göz Hastanesi Sincan",
        "Güven Hastanesi",
        "Medical Park (Yenimahalle)",
        "Medical Park Keçiören",
        "Mega Avrupa Diş Hastanesi",
        "TOBB ETÜ Hastanesi"
    ) }
)

$row = 1581
$firstNewRow = $row
foreach ($group in $groups) {
    foreach ($name in $group.Names) {
        $ws.Cells.Item($row, 1).Value = $il
        $ws.Cells.Item($row, 2).Value = $group.Tur
        $ws.Cells.Item($row, 3).Value = $name
        $row = $row + 1
    }
    # Blank separator row between sub-groups (matches existing layout style).
    $row = $row + 1
}

$lastNewRow = $row - 2

# Reflect the scrolled view / selection from the authored workbook.
$ws.Activate() | Out-Null
$ws.Range("A$($firstNewRow):C$($lastNewRow)").Select() | Out-Null
